$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 39, shifting existing rows 39:77 down to 40:78.
$ws.Rows(39).Insert()

# Populate the newly inserted row 39 with its data (same fixed columns as
# the rest of the table, plus the new record's own values).
$ws.Cells.Item(39, 1).Value = 1
$ws.Cells.Item(39, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(39, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(39, 4).Value = 44771
$ws.Cells.Item(39, 5).Value = 15
$ws.Cells.Item(39, 6).Value = 100114001
$ws.Cells.Item(39, 7).Value = "Papa"
$ws.Cells.Item(39, 8).Value = "Asterix"
$ws.Cells.Item(39, 9).Value = "1a (guarda)"
$ws.Cells.Item(39, 10).Value = 1000
$ws.Cells.Item(39, 11).Value = 10000
$ws.Cells.Item(39, 12).Value = 11000
$ws.Cells.Item(39, 13).Value = 10500
$ws.Cells.Item(39, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(39, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(39, 16).Value = 420
$ws.Cells.Item(39, 17).Value = 25
$ws.Cells.Item(39, 18).Value = "Hortaliza"
